$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 42.125
$ws.Range("I5").Value = 31.333334
$ws.Range("K5").Value = 31.333334
$ws.Range("M5").Value = 83.66666599999999

$ws.Range("H58").Value = 3970
$ws.Range("J58").Value = 6664.6
$ws.Range("L58").Value = 19993.8
$ws.Range("N58").Value = -20293.8

$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 30000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -31748

$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 90000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -98736

$ws.Range("H86").Value = 6929.0586
$ws.Range("J86").Value = 9015.1
$ws.Range("L86").Value = 9015.1
$ws.Range("N86").Value = -11261.1

$ws.Range("H89").Value = 6929.0586
$ws.Range("J89").Value = 9015.1
$ws.Range("L89").Value = 45075.5
$ws.Range("N89").Value = -56307.5

$ws.Range("H98").Value = 826.5484
$ws.Range("I98").Value = 826.5484
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 826.5484
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 671.4516
$ws.Range("N98").ClearContents()

$ws.Range("H112").Value = 334953.4
$ws.Range("J112").Value = 346472.47
$ws.Range("L112").Value = 1039417.41
$ws.Range("N112").Value = -1041633.41

$ws.Range("H115").Value = 726
$ws.Range("I115").Value = 621.64703
$ws.Range("K115").Value = 1864.94109
$ws.Range("M115").Value = -297.9410899999998

$ws.Range("H122").Value = 826.5484
$ws.Range("I122").Value = 826.5484
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2479.6452
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -29.64519999999993
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 867.05554
$ws.Range("I125").Value = 735.8
$ws.Range("J125").Value = 917.53845
$ws.Range("K125").Value = 6622.2
$ws.Range("L125").Value = 8257.84605
$ws.Range("M125").Value = -4162.2
$ws.Range("N125").Value = -13177.84605

$ws.Range("H137").Value = 434558.97
$ws.Range("J137").Value = 1192478.9
$ws.Range("L137").Value = 3577436.7
$ws.Range("N137").Value = -3582536.7

$ws.Range("H138").Value = 2357.6584
$ws.Range("I138").Value = 1692.8077
$ws.Range("J138").Value = 2666.3394
$ws.Range("K138").Value = 5078.4231
$ws.Range("L138").Value = 7999.0182
$ws.Range("M138").Value = 61.57690000000002
$ws.Range("N138").Value = -18279.0182

$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270

$ws.Range("H140").Value = 91990
$ws.Range("J140").Value = 91990
$ws.Range("L140").Value = 91990
$ws.Range("N140").Value = -102350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1650.7778
$ws.Range("J2").Value = 3633
$ws.Range("L2").Value = 3633
$ws.Range("N2").Value = -3859

$ws.Range("H32").Value = 5246.284
$ws.Range("I32").Value = 3398.1743
$ws.Range("J32").Value = 22906
$ws.Range("K32").Value = 3398.1743
$ws.Range("L32").Value = 22906
$ws.Range("M32").Value = -3111.1743
$ws.Range("N32").Value = -23480

$ws.Range("H43").Value = 25377
$ws.Range("J43").Value = 25377
$ws.Range("L43").Value = 25377
$ws.Range("N43").Value = -26003

$ws.Range("H45").Value = 14827.3125
$ws.Range("I45").Value = 16872.076
$ws.Range("K45").Value = 16872.076
$ws.Range("M45").Value = -16495.076

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H116").Value = 1650.7778
$ws.Range("J116").Value = 3633
$ws.Range("L116").Value = 3633
$ws.Range("N116").Value = -8221

$ws.Range("H122").Value = 3675.5833
$ws.Range("I122").Value = 3625.7693
$ws.Range("K122").Value = 10877.3079
$ws.Range("M122").Value = -8427.3079

$ws.Range("H132").Value = 2358.761
$ws.Range("I132").Value = 1932.0883
$ws.Range("K132").Value = 5796.2649
$ws.Range("M132").Value = -3266.2649

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1650.7778
$ws.Range("J3").Value = 3633
$ws.Range("L3").Value = 3633
$ws.Range("N3").Value = -3861

$ws.Range("H80").Value = 516.625
$ws.Range("I80").Value = 332
$ws.Range("J80").Value = 600.5454999999999
$ws.Range("K80").Value = 332
$ws.Range("L80").Value = 600.5454999999999
$ws.Range("M80").Value = 666
$ws.Range("N80").Value = -2596.5455

$ws.Range("H83").Value = 516.625
$ws.Range("I83").Value = 332
$ws.Range("J83").Value = 600.5454999999999
$ws.Range("K83").Value = 1660
$ws.Range("L83").Value = 3002.7275
$ws.Range("M83").Value = 3332
$ws.Range("N83").Value = -12986.7275

$ws.Range("H132").Value = 34701.703
$ws.Range("J132").Value = 34701.703
$ws.Range("L132").Value = 34701.703
$ws.Range("N132").Value = -44821.703

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 33519.22
$ws.Range("I134").Value = 2200.6086
$ws.Range("K134").Value = 6601.825800000001
$ws.Range("M134").Value = -4066.825800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6160.385
$ws.Range("I80").Value = 6693.5
$ws.Range("K80").Value = 6693.5
$ws.Range("M80").Value = -5695.5

$ws.Range("H83").Value = 6160.385
$ws.Range("I83").Value = 6693.5
$ws.Range("K83").Value = 33467.5
$ws.Range("M83").Value = -28475.5

$ws.Range("H122").Value = 12468.417
$ws.Range("I122").Value = 3559.1904
$ws.Range("K122").Value = 10677.5712
$ws.Range("M122").Value = -8227.5712

$ws.Range("H132").Value = 1919.1428
$ws.Range("I132").Value = 1763.6364
$ws.Range("J132").Value = 2489.3333
$ws.Range("K132").Value = 5290.9092
$ws.Range("L132").Value = 7467.999899999999
$ws.Range("M132").Value = -2760.9092
$ws.Range("N132").Value = -12527.9999

$ws.Range("H139").Value = 125000
$ws.Range("J139").Value = 125000
$ws.Range("L139").Value = 125000
$ws.Range("N139").Value = -135280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2734.276
$ws.Range("I132").Value = 2303.4
$ws.Range("K132").Value = 6910.200000000001
$ws.Range("M132").Value = -4380.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 13000
$ws.Range("J30").Value = 11333.333
$ws.Range("L30").Value = 11333.333
$ws.Range("N30").Value = -11547.333

$ws.Range("H81").Value = 2500.2
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 2500.2
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -40608

$ws.Range("H132").Value = 1674470.4
$ws.Range("I132").Value = 2118.3684
$ws.Range("J132").Value = 6213711.5
$ws.Range("K132").Value = 6355.1052
$ws.Range("L132").Value = 18641134.5
$ws.Range("M132").Value = -3825.1052
$ws.Range("N132").Value = -18646194.5
